# Auto-generated edit script: updates H:N leve-profit computed columns
# across multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW) to match
# refreshed market-board price data from the scheduled runner.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1146.8889
$ws.Range("I98").Value = 401.1579
$ws.Range("J98").Value = 2918
$ws.Range("K98").Value = 401.1579
$ws.Range("L98").Value = 2918
$ws.Range("M98").Value = 1096.8421
$ws.Range("N98").Value = -5914
$ws.Range("H100").Value = 30304460
$ws.Range("I100").Value = 1572
$ws.Range("J100").Value = 333333340
$ws.Range("K100").Value = 1572
$ws.Range("L100").Value = 333333340
$ws.Range("M100").Value = -1031
$ws.Range("N100").Value = -333334422
$ws.Range("H122").Value = 1146.8889
$ws.Range("I122").Value = 401.1579
$ws.Range("J122").Value = 2918
$ws.Range("K122").Value = 1203.4737
$ws.Range("L122").Value = 8754
$ws.Range("M122").Value = 1246.5263
$ws.Range("N122").Value = -13654
$ws.Range("H138").Value = 4406.1665
$ws.Range("I138").Value = 4597.5
$ws.Range("J138").Value = 4253.1
$ws.Range("K138").Value = 13792.5
$ws.Range("L138").Value = 12759.3
$ws.Range("M138").Value = -8652.5
$ws.Range("N138").Value = -23039.3

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22990.193
$ws.Range("I32").Value = 4387.2266
$ws.Range("K32").Value = 4387.2266
$ws.Range("M32").Value = -4100.2266
$ws.Range("H88").Value = 2669
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 2669
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 2669
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -3481
$ws.Range("H91").Value = 2669
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 2669
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 2669
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -5477
$ws.Range("H105").Value = 18566.666
$ws.Range("J105").Value = 18566.666
$ws.Range("L105").Value = 18566.666
$ws.Range("N105").Value = -25554.666
$ws.Range("H122").Value = 1514.15
$ws.Range("I122").Value = 1272.1875
$ws.Range("J122").Value = 2482
$ws.Range("K122").Value = 3816.5625
$ws.Range("L122").Value = 7446
$ws.Range("M122").Value = -1366.5625
$ws.Range("N122").Value = -12346

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2417.9412
$ws.Range("I86").Value = 2473.8667
$ws.Range("J86").Value = 1998.5
$ws.Range("K86").Value = 2473.8667
$ws.Range("L86").Value = 1998.5
$ws.Range("M86").Value = -1350.8667
$ws.Range("N86").Value = -4244.5
$ws.Range("H89").Value = 2417.9412
$ws.Range("I89").Value = 2473.8667
$ws.Range("J89").Value = 1998.5
$ws.Range("K89").Value = 12369.3335
$ws.Range("L89").Value = 9992.5
$ws.Range("M89").Value = -6753.333500000001
$ws.Range("N89").Value = -21224.5
$ws.Range("H105").Value = 3226.6296
$ws.Range("I105").Value = 1954.1177
$ws.Range("J105").Value = 5389.9
$ws.Range("K105").Value = 1954.1177
$ws.Range("L105").Value = 5389.9
$ws.Range("M105").Value = -207.1177
$ws.Range("N105").Value = -8883.9
$ws.Range("H107").Value = 8341.666999999999
$ws.Range("I107").Value = 12220
$ws.Range("K107").Value = 12220
$ws.Range("M107").Value = -10300

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 4177.421
$ws.Range("I105").Value = 3449.25
$ws.Range("J105").Value = 5425.7144
$ws.Range("K105").Value = 3449.25
$ws.Range("L105").Value = 5425.7144
$ws.Range("M105").Value = -1702.25
$ws.Range("N105").Value = -8919.714400000001
$ws.Range("H111").Value = 30560
$ws.Range("J111").Value = 30560
$ws.Range("L111").Value = 30560
$ws.Range("N111").Value = -38740

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 2235.7827
$ws.Range("I75").Value = 578.25
$ws.Range("J75").Value = 2584.7368
$ws.Range("K75").Value = 1734.75
$ws.Range("L75").Value = 7754.2104
$ws.Range("M75").Value = -736.75
$ws.Range("N75").Value = -9750.2104
$ws.Range("H78").Value = 2235.7827
$ws.Range("I78").Value = 578.25
$ws.Range("J78").Value = 2584.7368
$ws.Range("K78").Value = 5204.25
$ws.Range("L78").Value = 23262.6312
$ws.Range("M78").Value = -212.25
$ws.Range("N78").Value = -33246.6312
$ws.Range("H113").Value = 573.5106
$ws.Range("I113").Value = 538.4167
$ws.Range("J113").Value = 610.13043
$ws.Range("K113").Value = 1615.2501
$ws.Range("L113").Value = 1830.39129
$ws.Range("M113").Value = 554.7499
$ws.Range("N113").Value = -6170.39129
$ws.Range("H114").Value = 3565.2856
$ws.Range("I114").Value = 481
$ws.Range("J114").Value = 4799
$ws.Range("K114").Value = 1443
$ws.Range("L114").Value = 14397
$ws.Range("M114").Value = 1811
$ws.Range("N114").Value = -20905

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H45").Value = 14326.5
$ws.Range("J45").Value = 14326.5
$ws.Range("L45").Value = 14326.5
$ws.Range("N45").Value = -15444.5
$ws.Range("H51").Value = 14587.2
$ws.Range("J51").Value = 14587.2
$ws.Range("L51").Value = 14587.2
$ws.Range("N51").Value = -15605.2
$ws.Range("H113").Value = 1681.7273
$ws.Range("I113").Value = 1524.875
$ws.Range("K113").Value = 1524.875
$ws.Range("M113").Value = 645.125
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 678.3913
$ws.Range("I16").Value = 668.6818
$ws.Range("J16").Value = 892
$ws.Range("K16").Value = 668.6818
$ws.Range("L16").Value = 892
$ws.Range("M16").Value = -498.6818
$ws.Range("N16").Value = -1232
$ws.Range("H46").Value = 1097.3846
$ws.Range("I46").Value = 869.3333
$ws.Range("J46").Value = 1408.3636
$ws.Range("K46").Value = 869.3333
$ws.Range("L46").Value = 1408.3636
$ws.Range("M46").Value = -681.3333
$ws.Range("N46").Value = -1784.3636
